$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 267; this shifts existing rows 267-312 down to 268-313
$ws.Rows.Item(267).Insert()

# Populate the newly inserted row 267 with its data
$ws.Cells.Item(267, 1).Value = 11
$ws.Cells.Item(267, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(267, 3).Value = 'Bíobío'
$ws.Cells.Item(267, 4).Value = 45218
$ws.Cells.Item(267, 5).Value = 8
$ws.Cells.Item(267, 6).Value = 'Fruta'
$ws.Cells.Item(267, 7).Value = 100108
$ws.Cells.Item(267, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(267, 9).Value = 100108005
$ws.Cells.Item(267, 10).Value = 'Piña'
$ws.Cells.Item(267, 11).Value = 'Sin especificar'
$ws.Cells.Item(267, 12).Value = 'Segunda'
$ws.Cells.Item(267, 13).Value = 230
$ws.Cells.Item(267, 14).Value = 22000
$ws.Cells.Item(267, 15).Value = 22000
$ws.Cells.Item(267, 16).Value = 22000
$ws.Cells.Item(267, 17).Value = '$/caja 14 unidades'
$ws.Cells.Item(267, 18).Value = 'Ecuador'
$ws.Cells.Item(267, 19).Value = 1571
$ws.Cells.Item(267, 20).Value = 14

# Apply the same numeric style (date-time format) to D267 as the rest of column D
$ws.Cells.Item(267, 4).NumberFormat = $ws.Cells.Item(268, 4).NumberFormat
